$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selection to reflect the new range selected by the author.
[void]$ws.Range("B1:B128").Select()

# B1: previously an erroring formula; now a plain square curve (4*A1^2).
$ws.Range("B1").Formula = "=ROUND(4*A1*A1,0)"

# B2:B65 share one formula (si="0" in the OOXML); B66:B128 share another (si="1").
# Re-enter the formula across the whole B2:B128 range so Excel rebuilds the
# shared-formula groups exactly as the diff shows (two groups split at B66).
$ws.Range("B2:B65").Formula = "=ROUND(4*A2*A2,0)"
$ws.Range("B66:B128").Formula = "=ROUND(4*A66*A66,0)"

[void]$ws.Calculate()
